# Trade #111 closed at 2026-02-16 21:42:58 - leadlag DOWN +0.000%
#
# - Trade #87 (leadlag sheet, row 67) goes from OPEN -> CLOSED, and its
#   closed record is appended to the "All Trades" sheet.
# - Trade #111 is newly opened and appended to the "leadlag" sheet.
# - Summary / Comparison roll-up stats are refreshed to reflect the new
#   totals.

function Set-Text($ws, $ref, $val) {
    # Force a literal text value (Excel would otherwise auto-coerce
    # percent-looking / date-looking strings into numbers/dates).
    $r = $ws.Range($ref)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

function Set-Num($ws, $ref, $val) {
    $ws.Range($ref).Value = $val
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet: overall + leadlag roll-up numbers
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

Set-Num  $summary "C2" 87
Set-Text $summary "D2" "70.1%"
Set-Text $summary "E2" "+25.4685%"
Set-Text $summary "F2" "+0.2927%"

Set-Num  $summary "C3" 85
Set-Text $summary "D3" "49.4%"
Set-Text $summary "E3" "+14.1472%"
Set-Text $summary "F3" "+0.1664%"

# ---------------------------------------------------------------------
# leadlag sheet: close trade #87 (row 67) and append new trade #111
# ---------------------------------------------------------------------
$leadlag = $wb.Worksheets.Item("leadlag")

Set-Num  $leadlag "G67" 68460.857856
Set-Text $leadlag "H67" "CLOSED"
Set-Num  $leadlag "I67" 0.046
Set-Num  $leadlag "J67" 0.46
Set-Text $leadlag "M67" "time_exit_5min"
Set-Num  $leadlag "N67" 5

Set-Num  $leadlag "A87" 111
Set-Text $leadlag "B87" "2026-02-16"
Set-Text $leadlag "C87" "21:42:58"
Set-Text $leadlag "D87" "leadlag"
Set-Text $leadlag "E87" "DOWN"
Set-Num  $leadlag "F87" 68352.48
Set-Text $leadlag "H87" "OPEN"
Set-Num  $leadlag "I87" 0
Set-Num  $leadlag "J87" 0
Set-Num  $leadlag "K87" 0.6602
Set-Text $leadlag "L87" "Coinbase leading with -0.066% move"
Set-Num  $leadlag "N87" 0

# ---------------------------------------------------------------------
# All Trades sheet: append the now-closed trade #87
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

Set-Num  $allTrades "A88" 87
Set-Text $allTrades "B88" "2026-02-16"
Set-Text $allTrades "C88" "21:37:54"
Set-Text $allTrades "D88" "leadlag"
Set-Text $allTrades "E88" "DOWN"
Set-Num  $allTrades "F88" 68492.33500000001
Set-Num  $allTrades "G88" 68460.857856
Set-Text $allTrades "H88" "CLOSED"
Set-Num  $allTrades "I88" 0.046
Set-Num  $allTrades "J88" 0.46
Set-Num  $allTrades "K88" 0.7441
Set-Text $allTrades "L88" "Coinbase leading with -0.074% move"
Set-Text $allTrades "M88" "time_exit_5min"
Set-Num  $allTrades "N88" 5

# ---------------------------------------------------------------------
# Comparison sheet: leadlag strategy roll-up
# ---------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

Set-Num  $comparison "B2" 85
Set-Text $comparison "C2" "49.4%"
Set-Text $comparison "E2" "+0.5112%"
Set-Text $comparison "G2" "1.68"
